$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are plain-text strings that are NOT numeric-looking.
# These can be assigned directly; Excel will keep them as text.
$textUpdates = @{
    "D2" = '21.996.97'
    "E2" = '  -1.88%  '
    "D3" = '1.553.90'
    "E3" = '  -0.92%  '
    "E4" = '  -0.06%  '
    "E5" = '  -0.05%  '
    "E6" = '  -0.13%  '
    "E7" = '  +3.02%  '
    "E8" = '  -2.16%  '
    "E9" = '  -12.39%  '
    "E10" = '  -3.13%  '
    "E11" = '  -2.09%  '
    "E12" = '  -0.06%  '
    "E13" = '  -6.50%  '
    "E14" = '  -3.49%  '
    "E15" = '  -0.91%  '
    "D16" = '1.556.56'
    "E16" = '  -0.68%  '
    "E17" = '  -1.62%  '
    "E18" = '  -1.25%  '
    "E19" = '  -2.57%  '
    "E20" = '  +0.58%  '
    "E21" = '  -0.04%  '
    "E22" = '  -2.96%  '
    "E23" = '  -3.96%  '
    "D24" = '22.000.01'
    "E24" = '  -1.85%  '
    "E25" = '  -3.20%  '
    "E26" = '  -3.34%  '
    "E27" = '  -1.30%  '
    "E28" = '  -3.64%  '
    "E29" = '  -1.69%  '
    "D30" = '1.732.45'
    "E30" = '  -0.89%  '
    "E31" = '  -3.27%  '
    "E32" = '  +1.72%  '
    "E33" = '  -2.72%  '
    "B34" = 'FraxShare'
    "C34" = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
    "E34" = '  -5.33%  '
    "B35" = 'WEMIXTOKEN'
    "C35" = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
    "E35" = '  -16.77%  '
    "E36" = '  -1.96%  '
    "E37" = '  -2.61%  '
    "E38" = '  -0.51%  '
    "E40" = '  -4.49%  '
    "E41" = '  -5.67%  '
    "E42" = '  -4.20%  '
    "E43" = '  -0.09%  '
    "E44" = '  -3.63%  '
    "E45" = '  -2.38%  '
    "E46" = '  -1.10%  '
    "E47" = '  -4.07%  '
    "E48" = '  -5.09%  '
    "E49" = '  -4.34%  '
    "E50" = '  -3.27%  '
    "E51" = '  -4.30%  '
}

foreach ($ref in $textUpdates.Keys) {
    $ws.Range($ref).Value = $textUpdates[$ref]
}

# Cells whose new values look like plain numbers (e.g. "1.001", "13.50").
# These must be forced to Text number format first, otherwise Excel will
# auto-convert them to numeric values and silently drop formatting such as
# trailing zeros (e.g. "13.50" -> 13.5, "1.000" -> 1).
$numericLookingUpdates = @{
    "D5" = '1.001'
    "D6" = '286.58'
    "D7" = '0.3802'
    "D8" = '0.3238'
    "D9" = '41.43'
    "D10" = '1.120'
    "D11" = '0.07307'
    "D13" = '19.35'
    "D14" = '5.708'
    "D15" = '6.802'
    "D17" = '0.00001092'
    "D18" = '0.06625'
    "D19" = '85.12'
    "D20" = '6.411'
    "D21" = '0.9999'
    "D22" = '15.91'
    "D23" = '11.44'
    "D25" = '2.297'
    "D26" = '2.518'
    "D27" = '148.67'
    "D28" = '18.78'
    "D29" = '4.846'
    "D31" = '120.31'
    "D32" = '1.095'
    "D33" = '5.872'
    "D34" = '9.255'
    "D35" = '1.644'
    "D36" = '0.08134'
    "D37" = '0.06195'
    "D38" = '5.234'
    "D39" = '0.02284'
    "D40" = '0.2102'
    "D41" = '1.217'
    "D43" = '0.9999'
    "D44" = '0.5919'
    "D45" = '13.50'
    "D47" = '0.5733'
    "D49" = '119.24'
    "D50" = '1.155'
    "D51" = '0.06866'
}

foreach ($ref in $numericLookingUpdates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $numericLookingUpdates[$ref]
}
